$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 51, shifting existing rows 51-77 down to 52-78.
$ws.Rows.Item(51).Insert()

# Populate the newly inserted row 51 with the new weekly data point.
$ws.Range("A51").Value = 10
$ws.Range("B51").Value = "Vega Modelo de Temuco"
$ws.Range("C51").Value = "La Araucanía"
$ws.Range("D51").Value = 44830
$ws.Range("E51").Value = 9
$ws.Range("F51").Value = 100112026
$ws.Range("G51").Value = "Haba"
$ws.Range("H51").Value = "Sin especificar"
$ws.Range("I51").Value = "Primera"
$ws.Range("J51").Value = 65
$ws.Range("K51").Value = 12000
$ws.Range("L51").Value = 12000
$ws.Range("M51").Value = 12000
$ws.Range("N51").Value = "`$/saco 25 kilos"
$ws.Range("O51").Value = "Región Metropolitana"
$ws.Range("P51").Value = 480
$ws.Range("Q51").Value = 25
$ws.Range("R51").Value = "Hortaliza"
